$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fill in the weekly-hours table (B2:G8). Days with no hours logged get a
#    literal "-" (text), days with hours get the numeric value. The H column
#    SUM formulas and the row-9 totals recalc automatically.
# ---------------------------------------------------------------------------

# Andrea Favero
$ws.Range("B2").Value = "-"
$ws.Range("C2").Value = "-"
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 5
$ws.Range("G2").Value = 10

# Eleonora Thiella
$ws.Range("B3").Value = "-"
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "-"
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = 17

# Federico Caldart
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "-"
$ws.Range("F4").Value = "-"
$ws.Range("G4").Value = 13

# Giovanni Cavallin
$ws.Range("B5").Value = 8
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = "-"
$ws.Range("E5").Value = "-"
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 4

# Giovanni Dalla Riva
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 11

# Lorenzo Menegon
$ws.Range("B7").Value = "-"
$ws.Range("C7").Value = "-"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = "-"
$ws.Range("G7").Value = 18

# Stefano Panozzo
$ws.Range("B8").Value = "-"
$ws.Range("C8").Value = "-"
$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = "-"
$ws.Range("F8").Value = 4
$ws.Range("G8").Value = 17

# ---------------------------------------------------------------------------
# 2. Move / resize the chart (it now sits from col I row 1 to col V row 10
#    instead of col J/row1 to col R/row8). Compute the target Left/Top from
#    live column & row geometry so the anchor lands on the exact cell
#    boundaries, then add the EMU sub-cell offsets from the target anchor.
# ---------------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)

$fromLeft = $ws.Range("I1").Left + (428624 / 12700)
$fromTop  = $ws.Range("A1").Top  + (0 / 12700)
$toLeft   = $ws.Range("V1").Left + (285749 / 12700)
$toTop    = $ws.Range("A10").Top + (11906 / 12700)

$co.Left   = $fromLeft
$co.Top    = $fromTop
$co.Width  = $toLeft - $fromLeft
$co.Height = $toTop - $fromTop

# ---------------------------------------------------------------------------
# 3. The sheet no longer needs the stale external workbook link
#    (AnalisiRequisitiDettaglio.Orario.xlsx) - break it so the
#    externalReferences / externalLinks parts are dropped entirely.
# ---------------------------------------------------------------------------
$sources = $wb.LinkSources()
if ($sources) {
    foreach ($src in $sources) {
        $wb.BreakLink($src, 1)
    }
}

# ---------------------------------------------------------------------------
# 4. Leave the final selection on H11, matching where the user ended up
#    after finishing the data entry.
# ---------------------------------------------------------------------------
$ws.Range("H11").Select()
